$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the 2019 header block (W1:X1) into the new 2020
# header block (Y1:Z1), then set the year value.
$ws.Range("W1:X1").Copy()
$ws.Range("Y1:Z1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Y1").Value = 2020

# Copy the formatting of the 2019 data block (W2:X14) into the new 2020
# data block (Y2:Z14).
$ws.Range("W2:X14").Copy()
$ws.Range("Y2:Z14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the month labels (column Y, mirrors column A) and the count
# values (column Z) for the new 2020 block.
$months = "JAN","FEB","MAR","APR","MAY","JUN","JUL","AUG","SEPT","OCT","NOV","DEC"
for ($i = 0; $i -lt $months.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 25).Value = $months[$i]
    $ws.Cells.Item($row, 26).Value = -1
}

# Mirror the view state captured in the saved workbook: Z3:Z6 selected
# (as if the user had just highlighted the freshly-added 2020 counts).
$ws.Range("Z3:Z6").Select()
